$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty Advancement cells of the "Progressions" XML
# table with the literal text "none" (this also adds the new shared string
# and leaves the existing cell style (s="1") untouched).
$ws.Range("H2:K2").Value   = "none"
$ws.Range("D3:K3").Value   = "none"
$ws.Range("F4:K4").Value   = "none"
$ws.Range("D5:K5").Value   = "none"
$ws.Range("G6:K6").Value   = "none"
$ws.Range("D7:K7").Value   = "none"
$ws.Range("G8:K8").Value   = "none"
$ws.Range("D9:K9").Value   = "none"
$ws.Range("F10:K10").Value = "none"
$ws.Range("E11:K11").Value = "none"
$ws.Range("F12:K12").Value = "none"
$ws.Range("D13:K13").Value = "none"
$ws.Range("F14:K14").Value = "none"
$ws.Range("D15:K15").Value = "none"
$ws.Range("H16:K16").Value = "none"
$ws.Range("D17:K17").Value = "none"
$ws.Range("F18:K18").Value = "none"
$ws.Range("D19:K19").Value = "none"
$ws.Range("F20:K20").Value = "none"
$ws.Range("E21:K21").Value = "none"

# Move the active selection to where the author last left it.
$ws.Activate()
$ws.Range("H10").Select()
